# Add a new "2C3_Aluminum-production" sector row to the Sectors sheet.
# This mirrors the commit: "Added 2C_Aluminum-production sector."

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sectors")

# Insert a new row at position 44, pushing existing rows (old 44-66) down to 45-67.
$ws.Rows.Item(44).Insert()

# Populate the new row. Set column C (activity) before column A (sector) so the
# shared-string table allocates "aluminum_production" before
# "2C3_Aluminum-production", matching the canonical save order.
$ws.Range("C44").Value = "aluminum_production"
$ws.Range("A44").Value = "2C3_Aluminum-production"
$ws.Range("B44").Value = "process"
$ws.Range("D44").Value = "kt"
$ws.Range("E44").Value = "NC"

# Restore the view: selection moves to D45 (the cell that used to be D44,
# "1000", now shifted down one row by the insert).
$ws.Activate() | Out-Null
$ws.Range("D45").Select() | Out-Null
